# Added Jug of Wine run
# Appends a new pub-run record (row 57) to Sheet1, right below the last
# existing entry (row 56), mirroring its formatting, then selects the
# new row's first cell the way Excel leaves the selection after data entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the formatting (number formats / styles) of the previous data row
# (56) down onto the new row (57) for every column that will hold data,
# so the new cells pick up the same styles already used in the workbook
# (dates, times, plain text) instead of minting brand-new style entries.
foreach ($col in @("A","B","C","D","E","F","G","H","J","L","N","O","P")) {
    $ws.Range($col + "56").Copy()
    $ws.Range($col + "57").PasteSpecial(-4122)
}

# Fill in the new run's data.
$ws.Range("A57").Value = 43887
$ws.Range("B57").Value = "The Jug of Wine"
$ws.Range("C57").Value = "Lea"
$ws.Range("D57").Value = "start/end at pub"
$ws.Range("E57").Value = 3.82
$ws.Range("F57").Value = 0.042361111111111106
$ws.Range("G57").Formula = "=F57/E57"

# Re-apply G57's format after writing the formula: entering a formula that
# references the time-formatted F57 can otherwise pull in a time format.
$ws.Range("G56").Copy()
$ws.Range("G57").PasteSpecial(-4122)

$ws.Range("H57").Value = 1
$ws.Range("J57").Value = 1
$ws.Range("L57").Value = 1
$ws.Range("N57").Value = 1
$ws.Range("O57").Value = "Trespass, ""Oi, where you going?"", charming old church"
$ws.Range("P57").Formula = "=SUM(H57:N57)*E57"

# Leave the selection on the newly added row, as the author's session did.
$ws.Range("A57").Select()
